$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2014.8
$ws.Range("I137").Value = 1067.4667
$ws.Range("J137").Value = 3435.8
$ws.Range("K137").Value = 3202.4001
$ws.Range("L137").Value = 10307.4
$ws.Range("M137").Value = -652.4000999999998
$ws.Range("N137").Value = -15407.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1892
$ws.Range("I61").Value = 1892
$ws.Range("K61").Value = 1892
$ws.Range("M61").Value = -1680
$ws.Range("H74").Value = 4629.1665
$ws.Range("I74").Value = 1749.5
$ws.Range("J74").Value = 6069
$ws.Range("K74").Value = 1749.5
$ws.Range("L74").Value = 6069
$ws.Range("M74").Value = -875.5
$ws.Range("N74").Value = -7817
$ws.Range("H77").Value = 4629.1665
$ws.Range("I77").Value = 1749.5
$ws.Range("J77").Value = 6069
$ws.Range("K77").Value = 8747.5
$ws.Range("L77").Value = 30345
$ws.Range("M77").Value = -4379.5
$ws.Range("N77").Value = -39081
$ws.Range("H132").Value = 2521.4211
$ws.Range("I132").Value = 2521.4211
$ws.Range("K132").Value = 7564.263300000001
$ws.Range("M132").Value = -5034.263300000001
$ws.Range("H136").Value = 1892
$ws.Range("I136").Value = 1892
$ws.Range("K136").Value = 5676
$ws.Range("M136").Value = -3126
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 100000000
$ws.Range("I7").Value = 100000000
$ws.Range("K7").Value = 100000000
$ws.Range("M7").Value = -99999887
$ws.Range("H100").Value = 27306
$ws.Range("J100").Value = 27306
$ws.Range("L100").Value = 27306
$ws.Range("N100").Value = -29470
$ws.Range("H134").Value = 2423.75
$ws.Range("I134").Value = 2189.6365
$ws.Range("K134").Value = 6568.9095
$ws.Range("M134").Value = -4033.9095
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 9711
$ws.Range("J28").Value = 9614.666999999999
$ws.Range("L28").Value = 9614.666999999999
$ws.Range("N28").Value = -10104.667
$ws.Range("H141").Value = 19862.666
$ws.Range("J141").Value = 19862.666
$ws.Range("L141").Value = 19862.666
$ws.Range("N141").Value = -30222.666
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 3000
$ws.Range("K70").Value = 9000
$ws.Range("M70").Value = -8685
$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 3000
$ws.Range("K73").Value = 9000
$ws.Range("M73").Value = -7908
$ws.Range("H109").Value = 2799.6667
$ws.Range("J109").Value = 2999
$ws.Range("L109").Value = 8997
$ws.Range("N109").Value = -11077
$ws.Range("H131").Value = 1804.7174
$ws.Range("I131").Value = 1130
$ws.Range("J131").Value = 1835.3864
$ws.Range("K131").Value = 3390
$ws.Range("L131").Value = 5506.1592
$ws.Range("M131").Value = 1650
$ws.Range("N131").Value = -15586.1592
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 8899.4
$ws.Range("J36").Value = 7624.25
$ws.Range("L36").Value = 7624.25
$ws.Range("N36").Value = -8594.25
$ws.Range("H80").Value = 4421.125
$ws.Range("I80").Value = 4380.25
$ws.Range("J80").Value = 4462
$ws.Range("K80").Value = 4380.25
$ws.Range("L80").Value = 4462
$ws.Range("M80").Value = -3382.25
$ws.Range("N80").Value = -6458
$ws.Range("H83").Value = 4421.125
$ws.Range("I83").Value = 4380.25
$ws.Range("J83").Value = 4462
$ws.Range("K83").Value = 21901.25
$ws.Range("L83").Value = 22310
$ws.Range("M83").Value = -16909.25
$ws.Range("N83").Value = -32294
$ws.Range("H122").Value = 128659.625
$ws.Range("I122").Value = 3580.2
$ws.Range("J122").Value = 337125.34
$ws.Range("K122").Value = 10740.6
$ws.Range("L122").Value = 1011376.02
$ws.Range("M122").Value = -8290.599999999999
$ws.Range("N122").Value = -1016276.02
$ws.Range("H132").Value = 1980.1904
$ws.Range("I132").Value = 1634.8
$ws.Range("J132").Value = 8888
$ws.Range("K132").Value = 4904.4
$ws.Range("L132").Value = 26664
$ws.Range("M132").Value = -2374.4
$ws.Range("N132").Value = -31724
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4129.5386
$ws.Range("I22").Value = 3836.125
$ws.Range("J22").Value = 4599
$ws.Range("K22").Value = 3836.125
$ws.Range("L22").Value = 4599
$ws.Range("M22").Value = -3541.125
$ws.Range("N22").Value = -5189
$ws.Range("H27").Value = 4129.5386
$ws.Range("I27").Value = 3836.125
$ws.Range("J27").Value = 4599
$ws.Range("K27").Value = 3836.125
$ws.Range("L27").Value = 4599
$ws.Range("M27").Value = -3729.125
$ws.Range("N27").Value = -4813
$ws.Range("H43").Value = 1055555.5
$ws.Range("H74").Value = 43663.332
$ws.Range("I74").Value = 42995
$ws.Range("K74").Value = 42995
$ws.Range("M74").Value = -41997
$ws.Range("H77").Value = 43663.332
$ws.Range("I77").Value = 42995
$ws.Range("K77").Value = 128985
$ws.Range("M77").Value = -123993
$ws.Range("H132").Value = 3919.2903
$ws.Range("I132").Value = 3299.9167
$ws.Range("J132").Value = 6042.857
$ws.Range("K132").Value = 9899.750100000001
$ws.Range("L132").Value = 18128.571
$ws.Range("M132").Value = -7369.750100000001
$ws.Range("N132").Value = -23188.571
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7855.0557
$ws.Range("I62").Value = 5832
$ws.Range("J62").Value = 8259.666999999999
$ws.Range("K62").Value = 5832
$ws.Range("L62").Value = 8259.666999999999
$ws.Range("M62").Value = -5208
$ws.Range("N62").Value = -9507.666999999999
$ws.Range("H65").Value = 7855.0557
$ws.Range("I65").Value = 5832
$ws.Range("J65").Value = 8259.666999999999
$ws.Range("K65").Value = 29160
$ws.Range("L65").Value = 41298.335
$ws.Range("M65").Value = -26040
$ws.Range("N65").Value = -47538.335
$ws.Range("H98").Value = 45000
$ws.Range("J98").Value = 45000
$ws.Range("L98").Value = 45000
$ws.Range("N98").Value = -50990
$ws.Range("H132").Value = 2116
$ws.Range("I132").Value = 1827.5
$ws.Range("J132").Value = 2548.75
$ws.Range("K132").Value = 5482.5
$ws.Range("L132").Value = 7646.25
$ws.Range("M132").Value = -2952.5
$ws.Range("N132").Value = -12706.25
$ws.Range("H136").Value = 2476.5557
$ws.Range("I136").Value = 949
$ws.Range("J136").Value = 5073.4
$ws.Range("K136").Value = 2847
$ws.Range("L136").Value = 15220.2
$ws.Range("M136").Value = -297
$ws.Range("N136").Value = -20320.2
